$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the truncated email typo in-place (reuses the existing shared string).
$ws.Range("A2").Value = "sastaguvvu25@gmail.com"

# 2. Capture the current (post-fix) email values, top to bottom, before shifting
#    everything down a row to make room for a new header row.
$emails = @()
for ($r = 1; $r -le 5; $r++) {
    $emails += $ws.Cells.Item($r, 1).Value2
}

# 3. Hyperlinks.Delete() on this engine clears every hyperlink on the sheet, which
#    is exactly what we need since all five are being re-anchored one row down.
$ws.Range("A1").Hyperlinks.Delete()

# 4. Clear old contents, then lay out the new header + shifted data.
$ws.Range("A1:A5").ClearContents()
$ws.Range("A1").Value = "emails"
for ($i = 0; $i -lt $emails.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $emails[$i]
}

# 5. Restore styling: header is Normal (default), the five email cells are Hyperlink.
$ws.Range("A1").Style = "Normal"
$ws.Range("A2:A6").Style = "Hyperlink"

# 6. Re-add the hyperlinks at their new (shifted) positions, same target addresses,
#    in the same rId order as before the edit.
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:ravivarma25052@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:sampathsaicharan59@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:sahitya3066@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:sastaguvvu25@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:hfuafasfg@gmail.com") | Out-Null

# Re-apply the Hyperlink style again since Hyperlinks.Add() re-stamps its own xf.
$ws.Range("A2:A6").Style = "Hyperlink"

# 7. Fill the two new number columns.
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 2).Value = $r - 1
    $ws.Cells.Item($r, 3).Value = $r + 4
}

# 8. Match the saved selection state.
$ws.Range("E7").Select()
